$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme
for ($i=1; $i -le 8; $i++) {
    $c = $cs.Colors($i)
    $rgb = $c.RGB
    Write-Output "$i : $rgb  (hex: $('{0:X6}' -f $rgb))"
}
